$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D - shifts old D (Tên nhân viên) to E and old E (Ngày cấp) to F
$ws.Columns.Item(4).Insert()
# Work around an off-by-one in the engine's trailing default-column band after an insert
$ws.Columns.Item(16384).Delete()

# Set header text for the new column D2 (style already inherited as s4 from the insert)
$ws.Range("D2").Value = "Mã nhân viên"

# Column widths: new D column + widen old D (now E) column
$ws.Columns.Item(4).ColumnWidth = 18.6
$ws.Columns.Item(5).ColumnWidth = 35.0

# Update the saved selection to match the target state
$ws.Range("F8").Select()
